# 536-RBI-EPP-DB-SAR-REC-NON-RNI-CTRFD-DL-MD-TR-1-EarlyRePayment-Newcreateloan.xlsx
# Commit: "Loan RBI, Variable Instalments"
#
# On the "Repayment Schedule" sheet, insert a new (blank) column before the
# existing "Late" column, give it a custom width, and leave the sheet
# active with S8 selected (matching the author's last on-screen state).

$wb = $excel.ActiveWorkbook

$wsSchedule = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new, blank column at N -- pushes the old N ("Late") and
# O ("Outstanding") columns one to the right.
$wsSchedule.Columns("N:N").Insert()

# Give the freshly inserted column its own (non best-fit) width -- this
# rounds to the same pixel width Excel stores for the "In Advance" column
# (M), the closest representable value to the source width of 11.140625.
$wsSchedule.Columns("N:N").ColumnWidth = 10.33

# Leave "Repayment Schedule" as the active sheet/tab with S8 selected.
$wsSchedule.Activate() | Out-Null
$wsSchedule.Range("S8").Select() | Out-Null
